$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text (shared string "Puzzle 1" -> "Puzzle A")
$ws.Range("A1").Value = "Puzzle A"

# Green highlight color used for path-turn cells (RGB 0,176,80 -> FF00B050)
$greenColor = 5287936

# Fill in the solved puzzle path values (1-64) across C3:J10,
# highlighting the "turn" cells in green.

# Row 3
$ws.Range("C3").Value = 57
$ws.Range("C3").Interior.Color = $greenColor
$ws.Range("D3").Value = 58
$ws.Range("E3").Value = 59
$ws.Range("F3").Value = 60
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 62
$ws.Range("H3").Interior.Color = $greenColor
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 64
$ws.Range("J3").Interior.Color = $greenColor

# Row 4
$ws.Range("C4").Value = 56
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 32
$ws.Range("E4").Interior.Color = $greenColor
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 29
$ws.Range("I4").Value = 28
$ws.Range("J4").Value = 27

# Row 5
$ws.Range("C5").Value = 55
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = 1
$ws.Range("E5").Interior.Color = $greenColor
$ws.Range("F5").Value = 2
$ws.Range("F5").Interior.Color = $greenColor
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 22
$ws.Range("H5").Interior.Color = $greenColor
$ws.Range("I5").Value = 23
$ws.Range("J5").Value = 26

# Row 6
$ws.Range("C6").Value = 54
$ws.Range("D6").Value = 35
$ws.Range("D6").Interior.Color = $greenColor
$ws.Range("E6").Value = 36
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 19
$ws.Range("I6").Value = 24
$ws.Range("J6").Value = 25

# Row 7
$ws.Range("C7").Value = 53
$ws.Range("C7").Interior.Color = $greenColor
$ws.Range("D7").Value = 52
$ws.Range("E7").Value = 37
$ws.Range("F7").Value = 4
$ws.Range("F7").Interior.Color = $greenColor
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 18
$ws.Range("I7").Value = 17
$ws.Range("J7").Value = 16
$ws.Range("J7").Interior.Color = $greenColor

# Row 8
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 51
$ws.Range("E8").Value = 38
$ws.Range("F8").Value = 39
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 7
$ws.Range("I8").Value = 14
$ws.Range("J8").Value = 15

# Row 9
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 46
$ws.Range("E9").Value = 45
$ws.Range("F9").Value = 40
$ws.Range("G9").Value = 41
$ws.Range("H9").Value = 8
$ws.Range("H9").Interior.Color = $greenColor
$ws.Range("I9").Value = 13
$ws.Range("J9").Value = 12

# Row 10
$ws.Range("C10").Value = 48
$ws.Range("C10").Interior.Color = $greenColor
$ws.Range("D10").Value = 47
$ws.Range("E10").Value = 44
$ws.Range("E10").Interior.Color = $greenColor
$ws.Range("F10").Value = 43
$ws.Range("G10").Value = 42
$ws.Range("H10").Value = 9
$ws.Range("I10").Value = 10
$ws.Range("J10").Value = 11
$ws.Range("J10").Interior.Color = $greenColor
